$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price cells keep their exact text representation
# (Excel would otherwise reinterpret plain numeric-looking text as a number)
$textCells = @("D5","D6","D7","D10","D11","D12","D13","D14","D16","D17","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D37","D38","D39","D40","D41","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the data refresh
$ws.Range("D2").Value = "66.600.53"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "3.376.89"
$ws.Range("E3").Value = "  -4.77%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "561.84"
$ws.Range("E5").Value = "  -4.15%  "
$ws.Range("D6").Value = "183.91"
$ws.Range("E6").Value = "  -7.20%  "
$ws.Range("D7").Value = "0.601"
$ws.Range("E7").Value = "  -2.19%  "
$ws.Range("D9").Value = "3.369.11"
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").Value = "0.188"
$ws.Range("E10").Value = "  -8.86%  "
$ws.Range("D11").Value = "0.597"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").Value = "48.26"
$ws.Range("E12").Value = "  -7.47%  "
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  -7.07%  "
$ws.Range("D14").Value = "8.78"
$ws.Range("E14").Value = "  -6.04%  "
$ws.Range("D15").Value = "3.921.63"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "606.90"
$ws.Range("E16").Value = "  -11.73%  "
$ws.Range("D17").Value = "18.33"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "66.555.65"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("D19").Value = "3.383.16"
$ws.Range("E19").Value = "  -4.79%  "
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "11.63"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "0.920"
$ws.Range("E22").Value = "  -5.44%  "
$ws.Range("D23").Value = "17.00"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").Value = "5.23"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "98.53"
$ws.Range("E25").Value = "  -9.14%  "
$ws.Range("D26").Value = "4.10"
$ws.Range("E26").Value = "  -7.20%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "2.73"
$ws.Range("E28").Value = "  -7.69%  "
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -7.11%  "
$ws.Range("D30").Value = "8.84"
$ws.Range("E30").Value = "  -9.12%  "
$ws.Range("D31").Value = "30.84"
$ws.Range("E31").Value = "  -8.45%  "
$ws.Range("D32").Value = "3.88"
$ws.Range("E32").Value = "  -11.81%  "
$ws.Range("D33").Value = "6.36"
$ws.Range("E33").Value = "  -8.01%  "
$ws.Range("D34").Value = "11.20"
$ws.Range("E34").Value = "  -6.23%  "
$ws.Range("D35").Value = "555.09"
$ws.Range("E35").Value = "  +11.24%  "
$ws.Range("D36").Value = "3.883.62"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  -5.16%  "
$ws.Range("D38").Value = "58.20"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "3.45"
$ws.Range("E40").Value = "  -6.93%  "
$ws.Range("D41").Value = "3.51"
$ws.Range("E41").Value = "  +25.08%  "
$ws.Range("D42").Value = "0.0₃0725"
$ws.Range("E42").Value = "  -11.63%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  -8.00%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "0.129"
$ws.Range("E44").Value = "  -5.26%  "
$ws.Range("D45").Value = "0.351"
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("E46").Value = "  -6.15%  "
$ws.Range("D47").Value = "0.0421"
$ws.Range("E47").Value = "  -8.68%  "
$ws.Range("D48").Value = "3.22"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("D49").Value = "2.69"
$ws.Range("E49").Value = "  -8.58%  "
$ws.Range("D50").Value = "0.131"
$ws.Range("E50").Value = "  -4.57%  "
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.00%  "
